$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2024-01-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-06 Saturday", 2) | Out-Null

# Replace every arithmetic answer in the 20x5 practice table, cell by cell,
# so duplicate answer strings (e.g. "34+13=47") are each mapped to their own
# new value instead of being clobbered by a single global Find/Replace.
$answers = @(
    @("8+59=67", "1+97=98", "75-52=23", "5+11=16", "40-18=22"),
    @("83-44=39", "71-66=5", "51-48=3", "77-6=71", "76-27=49"),
    @("63+16=79", "44+33=77", "51+29=80", "58-29=29", "84-23=61"),
    @("63-21=42", "15+24=39", "99-0=99", "74-55=19", "26+61=87"),
    @("21+46=67", "47+21=68", "46+13=59", "65+25=90", "42-16=26"),
    @("11+57=68", "32+52=84", "12+35=47", "66-13=53", "90-20=70"),
    @("76-55=21", "30+66=96", "94-89=5", "60-45=15", "9+90=99"),
    @("31-4=27", "35+36=71", "82-9=73", "51+36=87", "41-26=15"),
    @("78-0=78", "62-24=38", "4+23=27", "74+13=87", "95-51=44"),
    @("23-3=20", "45+54=99", "0+58=58", "44+49=93", "98-31=67"),
    @("2+89=91", "8+84=92", "0+43=43", "24+34=58", "66-51=15"),
    @("34+16=50", "37+4=41", "82-14=68", "8+16=24", "79-16=63"),
    @("39+30=69", "71-33=38", "6+46=52", "75+20=95", "24+51=75"),
    @("99-33=66", "68-58=10", "48+47=95", "56+21=77", "54+0=54"),
    @("45+15=60", "21+50=71", "28-14=14", "43+17=60", "33+59=92"),
    @("63-2=61", "23-5=18", "4+19=23", "6+15=21", "69-11=58"),
    @("0+6=6", "64-26=38", "13+36=49", "29+5=34", "91-90=1"),
    @("36+21=57", "31+24=55", "33+15=48", "39+55=94", "19+53=72"),
    @("46-16=30", "83-5=78", "95+1=96", "93-83=10", "67-41=26"),
    @("86-39=47", "93-7=86", "44+14=58", "59-51=8", "38-22=16")
)

$t = $d.Tables.Item(1)
for ($r = 1; $r -le $answers.Length; $r++) {
    $row = $answers[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        $t.Cell($r, $c).Range.Text = $row[$c - 1]
    }
}

Write-Output "done"
